$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 689.8
$ws.Range("I28").Value = 488.77777
$ws.Range("J28").Value = 2499
$ws.Range("K28").Value = 488.77777
$ws.Range("L28").Value = 2499
$ws.Range("M28").Value = -3.777769999999975
$ws.Range("N28").Value = -3469

$ws.Range("H40").Value = 2137.4375
$ws.Range("J40").Value = 2230
$ws.Range("L40").Value = 2230
$ws.Range("N40").Value = -2580

$ws.Range("H74").Value = 5449.5
$ws.Range("I74").Value = 5449.5
$ws.Range("K74").Value = 5449.5
$ws.Range("M74").Value = -4513.5

$ws.Range("H77").Value = 5449.5
$ws.Range("I77").Value = 5449.5
$ws.Range("K77").Value = 27247.5
$ws.Range("M77").Value = -22567.5

$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()

$ws.Range("H113").Value = 1998.6
$ws.Range("I113").Value = 1999
$ws.Range("K113").Value = 1999
$ws.Range("M113").Value = 1255

$ws.Range("H141").Value = 3448.16
$ws.Range("I141").Value = 1873.1818
$ws.Range("K141").Value = 5619.5454
$ws.Range("M141").Value = -439.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3518.1638
$ws.Range("I32").Value = 2678.7856
$ws.Range("K32").Value = 2678.7856
$ws.Range("M32").Value = -2391.7856

$ws.Range("H63").Value = 2633.3333
$ws.Range("I63").Value = 2633.3333
$ws.Range("K63").Value = 2633.3333
$ws.Range("M63").Value = -1947.3333

$ws.Range("H66").Value = 2633.3333
$ws.Range("I66").Value = 2633.3333
$ws.Range("K66").Value = 13166.6665
$ws.Range("M66").Value = -9734.666499999999

$ws.Range("H74").Value = 494.02704
$ws.Range("I74").Value = 499.97223
$ws.Range("J74").Value = 280
$ws.Range("K74").Value = 499.97223
$ws.Range("L74").Value = 280
$ws.Range("M74").Value = 374.02777
$ws.Range("N74").Value = -2028

$ws.Range("H77").Value = 494.02704
$ws.Range("I77").Value = 499.97223
$ws.Range("J77").Value = 280
$ws.Range("K77").Value = 2499.86115
$ws.Range("L77").Value = 1400
$ws.Range("M77").Value = 1868.13885
$ws.Range("N77").Value = -10136

$ws.Range("H110").Value = 2062.647
$ws.Range("I110").Value = 2104.0625
$ws.Range("K110").Value = 2104.0625
$ws.Range("M110").Value = -59.0625

$ws.Range("H132").Value = 2493.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3644.2856
$ws.Range("I86").Value = 3424.3
$ws.Range("J86").Value = 4194.25
$ws.Range("K86").Value = 3424.3
$ws.Range("L86").Value = 4194.25
$ws.Range("M86").Value = -2301.3
$ws.Range("N86").Value = -6440.25

$ws.Range("H89").Value = 3644.2856
$ws.Range("I89").Value = 3424.3
$ws.Range("J89").Value = 4194.25
$ws.Range("K89").Value = 17121.5
$ws.Range("L89").Value = 20971.25
$ws.Range("M89").Value = -11505.5
$ws.Range("N89").Value = -32203.25

$ws.Range("H99").Value = 2162
$ws.Range("I99").Value = 1789.5555
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 1789.5555
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -291.5554999999999
$ws.Range("N99").Value = -5996

$ws.Range("H107").Value = 3298.8
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 11119.5
$ws.Range("J29").Value = 11119.5
$ws.Range("L29").Value = 11119.5
$ws.Range("N29").Value = -11705.5

$ws.Range("H31").Value = 2695.342
$ws.Range("I31").Value = 2152.3447
$ws.Range("K31").Value = 2152.3447
$ws.Range("M31").Value = -1857.3447

$ws.Range("H34").Value = 2695.342
$ws.Range("I34").Value = 2152.3447
$ws.Range("K34").Value = 2152.3447
$ws.Range("M34").Value = -1950.3447

$ws.Range("H45").Value = 1000
$ws.Range("I45").Value = 1000
$ws.Range("K45").Value = 1000
$ws.Range("M45").Value = -407

$ws.Range("H47").Value = 16023.333
$ws.Range("J47").Value = 15035
$ws.Range("L47").Value = 15035
$ws.Range("N47").Value = -16167

$ws.Range("H58").Value = 1709.0731
$ws.Range("I58").Value = 961.7273
$ws.Range("K58").Value = 961.7273
$ws.Range("M58").Value = -758.7273

$ws.Range("H132").Value = 2260
$ws.Range("I132").Value = 1964.5555
$ws.Range("K132").Value = 5893.666499999999
$ws.Range("M132").Value = -3363.666499999999

$ws.Range("H136").Value = 1709.0731
$ws.Range("I136").Value = 961.7273
$ws.Range("K136").Value = 2885.1819
$ws.Range("M136").Value = -335.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 4714.2856
$ws.Range("J62").Value = 3200
$ws.Range("L62").Value = 9600
$ws.Range("N62").Value = -10972

$ws.Range("H65").Value = 4714.2856
$ws.Range("J65").Value = 3200
$ws.Range("L65").Value = 28800
$ws.Range("N65").Value = -35664

$ws.Range("H113").Value = 1199
$ws.Range("I113").Value = 798.5
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 2395.5
$ws.Range("L113").Value = 6000
$ws.Range("M113").Value = -225.5
$ws.Range("N113").Value = -10340

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5768.4287
$ws.Range("J80").Value = 6665
$ws.Range("L80").Value = 6665
$ws.Range("N80").Value = -8661

$ws.Range("H83").Value = 5768.4287
$ws.Range("J83").Value = 6665
$ws.Range("L83").Value = 33325
$ws.Range("N83").Value = -43309

$ws.Range("H122").Value = 54373.42
$ws.Range("I122").Value = 1272.8
$ws.Range("J122").Value = 113374.11
$ws.Range("K122").Value = 3818.4
$ws.Range("L122").Value = 340122.33
$ws.Range("M122").Value = -1368.4
$ws.Range("N122").Value = -345022.33

$ws.Range("H126").Value = 2434.6667
$ws.Range("I126").Value = 2722.25
$ws.Range("J126").Value = 2204.6
$ws.Range("K126").Value = 8166.75
$ws.Range("L126").Value = 6613.799999999999
$ws.Range("M126").Value = -5696.75
$ws.Range("N126").Value = -11553.8

$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3169.75
$ws.Range("I132").Value = 2893
$ws.Range("K132").Value = 8679
$ws.Range("M132").Value = -6149

$ws.Range("H136").Value = 7774.6
$ws.Range("I136").Value = 7496.25
$ws.Range("K136").Value = 22488.75
$ws.Range("M136").Value = -19938.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1192.1818
$ws.Range("I136").Value = 748.8421
$ws.Range("K136").Value = 2246.5263
$ws.Range("M136").Value = 303.4737
